$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164069890975952
$ws.Range("B1").Value = 2.575162172317505
$ws.Range("C1").Value = 9.502006530761719
$ws.Range("D1").Value = 2.111528873443604
$ws.Range("E1").Value = 1.234969019889832
